# Efna1-Epha7.xlsx (NatmiData/natmiOut_TPM/YoungD4/LR-pairs_lrc2p)
# "update scripts wuth new tpm"
#
# The NATMI output was regenerated against a refreshed TPM matrix:
#   - the "Resolving-Mac" cluster used to appear only as a *target*
#     cluster (rows where column D = "Resolving-Mac"); it now appears
#     instead as a *sending* cluster (new rows 11-13, column A).
#   - every ligand/receptor expression, specificity, and edge-weight
#     statistic (columns G-T) is recomputed for rows 2-13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2: ECs -> ECs (Efna1/Epha7)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 30.194115
$ws.Range("H2").Value = 90.582345
$ws.Range("I2").Value = 0.9018420607989291
$ws.Range("J2").Value = 0.901842060798929
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2282403333333334
$ws.Range("N2").Value = 0.684721
$ws.Range("O2").Value = 0.1514399067192403
$ws.Range("P2").Value = 0.1514399067192403
$ws.Range("Q2").Value = 6.891514872305001
$ws.Range("R2").Value = 62.02363385074501
$ws.Range("S2").Value = 0.1365748775628773
$ws.Range("T2").Value = 0.1365748775628773

# row 3: ECs -> FAPs (Efna1/Epha7)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 30.194115
$ws.Range("H3").Value = 90.582345
$ws.Range("I3").Value = 0.9018420607989291
$ws.Range("J3").Value = 0.901842060798929
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.151356333333333
$ws.Range("N3").Value = 3.454069
$ws.Range("O3").Value = 0.763937263734893
$ws.Range("P3").Value = 0.763937263734893
$ws.Range("Q3").Value = 34.764185534645
$ws.Range("R3").Value = 312.877669811805
$ws.Range("S3").Value = 0.6889507562477709
$ws.Range("T3").Value = 0.6889507562477709

# row 4: ECs -> MuSCs (Efna1/Epha7)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha7"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 30.194115
$ws.Range("H4").Value = 90.582345
$ws.Range("I4").Value = 0.9018420607989291
$ws.Range("J4").Value = 0.901842060798929
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.127538
$ws.Range("N4").Value = 0.3826140000000001
$ws.Range("O4").Value = 0.08462282954586674
$ws.Range("P4").Value = 0.08462282954586674
$ws.Range("Q4").Value = 3.85089703887
$ws.Range("R4").Value = 34.65807334983
$ws.Range("S4").Value = 0.07631642698828096
$ws.Range("T4").Value = 0.07631642698828096

# row 5: FAPs -> ECs (Efna1/Epha7)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.633202
$ws.Range("H5").Value = 7.899606
$ws.Range("I5").Value = 0.07864884657754871
$ws.Range("J5").Value = 0.07864884657754868
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.2282403333333334
$ws.Range("N5").Value = 0.684721
$ws.Range("O5").Value = 0.1514399067192403
$ws.Range("P5").Value = 0.1514399067192403
$ws.Range("Q5").Value = 0.6010029022140001
$ws.Range("R5").Value = 5.409026119926001
$ws.Range("S5").Value = 0.01191057398927982
$ws.Range("T5").Value = 0.01191057398927981

# row 6: FAPs -> FAPs (Efna1/Epha7)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.633202
$ws.Range("H6").Value = 7.899606
$ws.Range("I6").Value = 0.07864884657754871
$ws.Range("J6").Value = 0.07864884657754868
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.151356333333333
$ws.Range("N6").Value = 3.454069
$ws.Range("O6").Value = 0.763937263734893
$ws.Range("P6").Value = 0.763937263734893
$ws.Range("Q6").Value = 3.031753799646001
$ws.Range("R6").Value = 27.28578419681401
$ws.Range("S6").Value = 0.06008278465035797
$ws.Range("T6").Value = 0.06008278465035795

# row 7: FAPs -> MuSCs (Efna1/Epha7)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha7"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.633202
$ws.Range("H7").Value = 7.899606
$ws.Range("I7").Value = 0.07864884657754871
$ws.Range("J7").Value = 0.07864884657754868
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.127538
$ws.Range("N7").Value = 0.3826140000000001
$ws.Range("O7").Value = 0.08462282954586674
$ws.Range("P7").Value = 0.08462282954586674
$ws.Range("Q7").Value = 0.3358333166760001
$ws.Range("R7").Value = 3.022499850084001
$ws.Range("S7").Value = 0.006655487937910929
$ws.Range("T7").Value = 0.006655487937910927

# row 8: MuSCs -> ECs (Efna1/Epha7)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5479553333333333
$ws.Range("H8").Value = 1.643866
$ws.Range("I8").Value = 0.01636640673320273
$ws.Range("J8").Value = 0.01636640673320272
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.2282403333333334
$ws.Range("N8").Value = 0.684721
$ws.Range("O8").Value = 0.1514399067192403
$ws.Range("P8").Value = 0.1514399067192403
$ws.Range("Q8").Value = 0.1250655079317778
$ws.Range("R8").Value = 1.125589571386
$ws.Range("S8").Value = 0.002478527109005368
$ws.Range("T8").Value = 0.002478527109005367

# row 9: MuSCs -> FAPs (Efna1/Epha7)
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5479553333333333
$ws.Range("H9").Value = 1.643866
$ws.Range("I9").Value = 0.01636640673320273
$ws.Range("J9").Value = 0.01636640673320272
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.151356333333333
$ws.Range("N9").Value = 3.454069
$ws.Range("O9").Value = 0.763937263734893
$ws.Range("P9").Value = 0.763937263734893
$ws.Range("Q9").Value = 0.6308918434171111
$ws.Range("R9").Value = 5.678026590754001
$ws.Range("S9").Value = 0.01250290797693522
$ws.Range("T9").Value = 0.01250290797693522

# row 10: MuSCs -> MuSCs (Efna1/Epha7)
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha7"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5479553333333333
$ws.Range("H10").Value = 1.643866
$ws.Range("I10").Value = 0.01636640673320273
$ws.Range("J10").Value = 0.01636640673320272
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.127538
$ws.Range("N10").Value = 0.3826140000000001
$ws.Range("O10").Value = 0.08462282954586674
$ws.Range("P10").Value = 0.08462282954586674
$ws.Range("Q10").Value = 0.06988512730266667
$ws.Range("R10").Value = 0.6289661457240001
$ws.Range("S10").Value = 0.00138497164726214
$ws.Range("T10").Value = 0.00138497164726214

# row 11: Resolving-Mac -> ECs (Efna1/Epha7)
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha7"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1052186666666667
$ws.Range("H11").Value = 0.315656
$ws.Range("I11").Value = 0.00314268589031943
$ws.Range("J11").Value = 0.003142685890319429
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.2282403333333334
$ws.Range("N11").Value = 0.684721
$ws.Range("O11").Value = 0.1514399067192403
$ws.Range("P11").Value = 0.1514399067192403
$ws.Range("Q11").Value = 0.02401514355288889
$ws.Range("R11").Value = 0.216136291976
$ws.Range("S11").Value = 0.0004759280580778472
$ws.Range("T11").Value = 0.000475928058077847

# row 12: Resolving-Mac -> FAPs (Efna1/Epha7)
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha7"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1052186666666667
$ws.Range("H12").Value = 0.315656
$ws.Range("I12").Value = 0.00314268589031943
$ws.Range("J12").Value = 0.003142685890319429
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.151356333333333
$ws.Range("N12").Value = 3.454069
$ws.Range("O12").Value = 0.763937263734893
$ws.Range("P12").Value = 0.763937263734893
$ws.Range("Q12").Value = 0.1211441782515556
$ws.Range("R12").Value = 1.090297604264
$ws.Range("S12").Value = 0.002400814859828882
$ws.Range("T12").Value = 0.002400814859828881

# row 13: Resolving-Mac -> MuSCs (Efna1/Epha7)
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha7"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1052186666666667
$ws.Range("H13").Value = 0.315656
$ws.Range("I13").Value = 0.00314268589031943
$ws.Range("J13").Value = 0.003142685890319429
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.127538
$ws.Range("N13").Value = 0.3826140000000001
$ws.Range("O13").Value = 0.08462282954586674
$ws.Range("P13").Value = 0.08462282954586674
$ws.Range("Q13").Value = 0.01341937830933333
$ws.Range("R13").Value = 0.120774404784
$ws.Range("S13").Value = 0.0002659429724127016
$ws.Range("T13").Value = 0.0002659429724127015

